$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# Replace "konstelasi grus" with "rasi bintang Grus" everywhere in the
# document body (matches case so we don't clobber capitalization rules).
$find.Execute("konstelasi grus", $true, $true, $false, $false, $false, `
               $true, 1, $false, "rasi bintang Grus", 2)
